$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in A1 ("ShortPName") is renamed to reflect how the IS is
# literally named in the contract.
$ws.Range("A1").Value = "ID_DES.LM_project"

# Column A needs to be a bit wider to fit the new header text.
$ws.Columns.Item(1).ColumnWidth = 18.7

# Leave the selection on A2, as the author left it.
[void]$ws.Range("A2").Select()
